$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    # Force the cell to Text format before assignment so numeric-looking
    # strings (e.g. "0.120", "0.0000158", "65.463.40") are stored verbatim
    # instead of being coerced into floating point / scientific notation.
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '65.463.40'
$ws.Range('E2').Value = '  +1.04%  '
Set-TextValue 'D3' '3.196.48'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue 'D5' '576.84'
$ws.Range('E5').Value = '  -0.18%  '
Set-TextValue 'D6' '168.20'
$ws.Range('E6').Value = '  -2.55%  '
Set-TextValue 'D7' '0.598'
$ws.Range('E7').Value = '  -5.19%  '
$ws.Range('E8').Value = '  -0.01%  '
Set-TextValue 'D9' '0.120'
$ws.Range('E9').Value = '  -2.01%  '
Set-TextValue 'D10' '6.72'
$ws.Range('E10').Value = '  -0.70%  '
Set-TextValue 'D11' '0.389'
$ws.Range('E11').Value = '  +0.19%  '
Set-TextValue 'D12' '3.753.63'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  -0.31%  '
Set-TextValue 'D14' '65.391.79'
$ws.Range('E14').Value = '  +0.78%  '
Set-TextValue 'D15' '25.68'
$ws.Range('E15').Value = '  +0.01%  '
Set-TextValue 'D16' '3.194.49'
$ws.Range('E16').Value = '  -1.27%  '
Set-TextValue 'D17' '0.0000158'
$ws.Range('E17').Value = '  -0.63%  '
Set-TextValue 'D18' '413.83'
$ws.Range('E18').Value = '  -0.75%  '
Set-TextValue 'D19' '12.86'
$ws.Range('E19').Value = '  +0.37%  '
Set-TextValue 'D20' '5.34'
$ws.Range('E20').Value = '  -0.65%  '
Set-TextValue 'D21' '7.17'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('E22').Value = '  +0.12%  '
Set-TextValue 'D23' '69.37'
$ws.Range('E23').Value = '  -2.12%  '
Set-TextValue 'D24' '0.202'
$ws.Range('E24').Value = '  -1.95%  '
Set-TextValue 'D25' '0.490'
$ws.Range('E25').Value = '  -0.99%  '
Set-TextValue 'D26' '0.0000106'
$ws.Range('E26').Value = '  -4.33%  '
Set-TextValue 'D27' '8.96'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('E28').Value = '  -0.06%  '
Set-TextValue 'D29' '1.84'
$ws.Range('E29').Value = '  -1.55%  '
Set-TextValue 'D30' '21.56'
$ws.Range('E30').Value = '  -1.32%  '
Set-TextValue 'D31' '5.03'
$ws.Range('E31').Value = '  +0.70%  '
Set-TextValue 'D32' '6.40'
$ws.Range('E32').Value = '  -0.21%  '
Set-TextValue 'D33' '1.15'
$ws.Range('E33').Value = '  -1.04%  '
Set-TextValue 'D34' '156.52'
$ws.Range('E34').Value = '  -0.65%  '
Set-TextValue 'D35' '1.37'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('B36').Value = 'Stacks'
$ws.Range('C36').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D36' '1.74'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D37' '2.737.96'
$ws.Range('E37').Value = '  -3.41%  '
Set-TextValue 'D38' '24.24'
$ws.Range('E38').Value = '  -4.24%  '
Set-TextValue 'D39' '4.17'
$ws.Range('E39').Value = '  -1.72%  '
Set-TextValue 'D40' '0.711'
$ws.Range('E40').Value = '  -1.98%  '
Set-TextValue 'D41' '0.0636'
$ws.Range('E41').Value = '  +0.98%  '
Set-TextValue 'D42' '5.58'
$ws.Range('E42').Value = '  -2.69%  '
Set-TextValue 'D43' '0.0264'
$ws.Range('E43').Value = '  +0.21%  '
Set-TextValue 'D44' '295.85'
$ws.Range('E44').Value = '  -1.70%  '
Set-TextValue 'D45' '21.64'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('E46').Value = '  -0.03%  '
Set-TextValue 'D47' '0.0991'
$ws.Range('E47').Value = '  -2.11%  '
Set-TextValue 'D48' '1.99'
$ws.Range('E48').Value = '  -8.10%  '
Set-TextValue 'D49' '5.81'
$ws.Range('E49').Value = '  -0.34%  '
Set-TextValue 'D50' '10.48'
$ws.Range('E50').Value = '  +0.80%  '
Set-TextValue 'D51' '0.906'
$ws.Range('E51').Value = '  -2.55%  '
